# ItemBasicMulStat 데이터 float으로 변경
# Change the declared type of columns C-K on the Item_BasicMulStat row 2
# (type-descriptor row) to "float", and update the saved selections that
# Excel records for the Item_BasicMulStat and ForReference_Storage sheets.

$wb = $excel.ActiveWorkbook

$wsMul = $wb.Worksheets.Item("Item_BasicMulStat")
$wsMul.Activate()

# Row 2 holds the per-column data type label. Every column except B
# (already "float") moves from "int"/"long long" to "float".
$wsMul.Range("C2:K2").Value = "float"

# Update the remembered selection for this sheet.
$wsMul.Range("J10").Select()

$wsStorage = $wb.Worksheets.Item("ForReference_Storage")
$wsStorage.Activate()
$wsStorage.Range("C34").Select()

# Restore focus back to the originally active sheet.
$wsMul.Activate()
